$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ G = 1.02; M = 0.63; S = 0.5 }
    3  = @{ G = 1.02; M = 0.63; S = 0.5 }
    4  = @{ G = 0.77; M = 0.53; S = 0.48 }
    5  = @{ G = 0.77; M = 0.53; S = 0.48 }
    6  = @{ G = 0.77; M = 0.53; S = 0.47 }
    7  = @{ G = 0.66; M = 0.47; S = 0.18 }
    8  = @{ G = 0.66; M = 0.47; S = 0.18 }
    9  = @{ G = 0.77; M = 0.53; S = 0.47 }
    10 = @{ G = 0.77; M = 0.53; S = 0.48 }
    11 = @{ G = 0.77; M = 0.53; S = 0.48 }
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Range("G$row").Value = $rowVals.G
    $ws.Range("M$row").Value = $rowVals.M
    $ws.Range("S$row").Value = $rowVals.S
}
